$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 223, shifting existing rows 223:315 down to 224:316.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with a new Ciboulette price record.
$ws.Range("A223").Value = 3
$ws.Range("B223").Value = "Femacal de La Calera"
$ws.Range("C223").Value = "Coquimbo"
$ws.Range("D223").Value = 44704
$ws.Range("E223").Value = 5
$ws.Range("F223").Value = 100112039
$ws.Range("G223").Value = "Ciboulette"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 120
$ws.Range("K223").Value = 1500
$ws.Range("L223").Value = 1500
$ws.Range("M223").Value = 1500
$ws.Range("N223").Value = "`$/docena de atados"
$ws.Range("O223").Value = "Provincia de Quillota"
$ws.Range("P223").Value = 500
$ws.Range("Q223").Value = 3
$ws.Range("R223").Value = "Hortaliza"
